$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1 (20:50 -> 21:20)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 21:20"

# Update Murcia's row (row 29) statistics
$ws.Range("B29").Value = 1084
$ws.Range("C29").Value = 45
$ws.Range("D29").Value = 999
$ws.Range("E29").Value = 42
